$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-PlainText 2 4 '24.950.55'
Set-PlainText 2 5 '  -3.97%  '

Set-PlainText 3 4 '1.642.21'
Set-PlainText 3 5 '  -5.63%  '

Set-TextValue 4 4 '0.9994'
Set-PlainText 4 5 '  -0.15%  '

Set-TextValue 5 4 '232.84'
Set-PlainText 5 5 '  -5.43%  '

Set-PlainText 6 5 '  -0.08%  '

Set-TextValue 7 4 '0.4737'
Set-PlainText 7 5 '  -5.66%  '

Set-PlainText 8 2 'OKB'
Set-PlainText 8 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 8 4 '39.41'
Set-PlainText 8 5 '  -3.56%  '

Set-PlainText 9 2 'Cardano'
Set-PlainText 9 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 9 4 '0.2580'
Set-PlainText 9 5 '  -5.71%  '

Set-PlainText 10 2 'Dogecoin'
Set-PlainText 10 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 10 4 '0.06096'
Set-PlainText 10 5 '  -1.38%  '

Set-PlainText 11 2 'TRON'
Set-PlainText 11 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 11 4 '0.07030'
Set-PlainText 11 5 '  -3.05%  '

Set-PlainText 12 2 'WrappedEther'
Set-PlainText 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-PlainText 12 4 '1.647.00'
Set-PlainText 12 5 '  -5.39%  '

Set-PlainText 13 2 'Solana'
Set-PlainText 13 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 13 4 '14.49'
Set-PlainText 13 5 '  -3.88%  '

Set-PlainText 14 2 'Polygon'
Set-PlainText 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 14 4 '0.5868'
Set-PlainText 14 5 '  -10.14%  '

Set-PlainText 15 2 'Polkadot'
Set-PlainText 15 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 15 4 '4.335'
Set-PlainText 15 5 '  -7.35%  '

Set-PlainText 16 2 'Litecoin'
Set-PlainText 16 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 16 4 '73.60'
Set-PlainText 16 5 '  -5.13%  '

Set-PlainText 17 2 'Dai'
Set-PlainText 17 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 17 4 '1.000'
Set-PlainText 17 5 '  -0.07%  '

Set-PlainText 18 2 'BinanceUSD'
Set-PlainText 18 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 18 4 '1.001'
Set-PlainText 18 5 '  -0.03%  '

Set-PlainText 19 2 'WrappedBTC'
Set-PlainText 19 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-PlainText 19 4 '24.958.21'
Set-PlainText 19 5 '  -4.07%  '

Set-PlainText 20 2 'ShibaInu'
Set-PlainText 20 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 20 4 '0.000006580'
Set-PlainText 20 5 '  -3.96%  '

Set-PlainText 21 2 'Avalanche'
Set-PlainText 21 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 21 4 '11.20'
Set-PlainText 21 5 '  -6.05%  '

Set-PlainText 22 2 'WrappedliquidstakedEther2.0'
Set-PlainText 22 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-PlainText 22 4 '1.857.63'
Set-PlainText 22 5 '  -5.70%  '

Set-PlainText 23 2 'Uniswap'
Set-PlainText 23 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 23 4 '4.285'
Set-PlainText 23 5 '  -4.72%  '

Set-PlainText 24 2 'Cosmos'
Set-PlainText 24 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 24 4 '8.548'
Set-PlainText 24 5 '  -1.74%  '

Set-PlainText 25 2 'Chainlink'
Set-PlainText 25 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 25 4 '5.225'
Set-PlainText 25 5 '  -3.17%  '

Set-PlainText 26 2 'Monero'
Set-PlainText 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 26 4 '133.95'
Set-PlainText 26 5 '  -1.15%  '

Set-PlainText 27 2 'EthereumClassic'
Set-PlainText 27 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 27 4 '14.88'
Set-PlainText 27 5 '  -2.57%  '

Set-PlainText 28 2 'Toncoin'
Set-PlainText 28 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 28 4 '1.383'
Set-PlainText 28 5 '  -8.59%  '

Set-PlainText 29 2 'BitcoinCash'
Set-PlainText 29 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 29 4 '103.21'
Set-PlainText 29 5 '  -2.13%  '

Set-PlainText 30 2 'LidoDAOToken'
Set-PlainText 30 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 30 4 '1.630'
Set-PlainText 30 5 '  -8.56%  '

Set-PlainText 31 2 'InternetComputer(DFINITY)'
Set-PlainText 31 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 31 4 '3.874'
Set-PlainText 31 5 '  -1.91%  '

Set-PlainText 32 2 'Filecoin'
Set-PlainText 32 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 32 4 '3.573'
Set-PlainText 32 5 '  -2.77%  '

Set-PlainText 33 2 'Stellar'
Set-PlainText 33 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 33 4 '0.07584'
Set-PlainText 33 5 '  -6.92%  '

Set-PlainText 34 2 'Frax'
Set-PlainText 34 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 34 4 '0.9995'
Set-PlainText 34 5 '  -0.05%  '

Set-PlainText 35 2 'Hedera'
Set-PlainText 35 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 35 4 '0.04273'
Set-PlainText 35 5 '  -9.05%  '

Set-PlainText 36 2 'HuobiToken'
Set-PlainText 36 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 36 4 '2.573'
Set-PlainText 36 5 '  -3.52%  '

Set-PlainText 37 2 'ARBITRUM'
Set-PlainText 37 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 37 4 '0.9248'
Set-PlainText 37 5 '  -6.91%  '

Set-PlainText 38 2 'ImmutableX'
Set-PlainText 38 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 38 4 '0.5892'
Set-PlainText 38 5 '  -3.34%  '

Set-PlainText 39 2 'MXToken'
Set-PlainText 39 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 39 4 '2.582'
Set-PlainText 39 5 '  -6.58%  '

Set-PlainText 40 2 'TrustWalletToken'
Set-PlainText 40 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 40 4 '0.8699'
Set-PlainText 40 5 '  +9.48%  '

Set-PlainText 41 2 'PaxDollar'
Set-PlainText 41 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 41 4 '1.0000'
Set-PlainText 41 5 '  -0.10%  '

Set-PlainText 42 2 'VeChain'
Set-PlainText 42 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 42 4 '0.01501'
Set-PlainText 42 5 '  -7.24%  '

Set-PlainText 43 2 'Quant'
Set-PlainText 43 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 43 4 '98.65'
Set-PlainText 43 5 '  -2.33%  '

Set-PlainText 44 2 'RenderToken'
Set-PlainText 44 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 44 4 '1.755'
Set-PlainText 44 5 '  -8.75%  '

Set-PlainText 45 2 'TheSandbox'
Set-PlainText 45 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 45 4 '0.3692'
Set-PlainText 45 5 '  -5.27%  '

Set-PlainText 46 2 'FraxShare'
Set-PlainText 46 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 46 4 '4.671'
Set-PlainText 46 5 '  -6.81%  '

Set-PlainText 47 2 'Algorand'
Set-PlainText 47 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 47 4 '0.1101'
Set-PlainText 47 5 '  -5.46%  '

Set-PlainText 48 2 'Aptos'
Set-PlainText 48 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 48 4 '6.089'
Set-PlainText 48 5 '  -3.47%  '

Set-PlainText 49 2 'Cronos'
Set-PlainText 49 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 49 4 '0.05205'
Set-PlainText 49 5 '  -1.67%  '

Set-PlainText 51 2 'Elrond'
Set-PlainText 51 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 51 4 '28.70'
Set-PlainText 51 5 '  -6.82%  '
